$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notes = $s.NotesPage
$ph = $notes.Shapes.AddPlaceholder(2)
$ph.TextFrame.TextRange.Text = "This is a test"
